# Refresh cached Universalis market-board figures (currentAveragePrice*, LevePrice*,
# LeveProfit*) for the leve-crafting profit tracker across all eight job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 80
$ws.Range("H80").Value = 1387.5714
$ws.Range("I80").Value = 1166
$ws.Range("J80").Value = 1553.75
$ws.Range("K80").Value = 3498
$ws.Range("L80").Value = 4661.25
$ws.Range("M80").Value = -2500
$ws.Range("N80").Value = -6657.25
# Row 83
$ws.Range("H83").Value = 1387.5714
$ws.Range("I83").Value = 1166
$ws.Range("J83").Value = 1553.75
$ws.Range("K83").Value = 10494
$ws.Range("L83").Value = 13983.75
$ws.Range("M83").Value = -5502
$ws.Range("N83").Value = -23967.75
# Row 87
$ws.Range("H87").Value = 88399.5
$ws.Range("J87").Value = 88399.5
$ws.Range("L87").Value = 88399.5
$ws.Range("N87").Value = -90895.5
# Row 90
$ws.Range("H90").Value = 88399.5
$ws.Range("J90").Value = 88399.5
$ws.Range("L90").Value = 265198.5
$ws.Range("N90").Value = -277678.5
# Row 120
$ws.Range("H120").Value = 62480
$ws.Range("J120").Value = 62480
$ws.Range("L120").Value = 62480
$ws.Range("N120").Value = -72156
# Row 132
$ws.Range("H132").Value = 1796.6923
$ws.Range("I132").Value = 1796.6923
$ws.Range("K132").Value = 5390.0769
$ws.Range("M132").Value = -2860.0769
# Row 138
$ws.Range("H138").Value = 1957.7677
$ws.Range("J138").Value = 2186.9
$ws.Range("L138").Value = 6560.700000000001
$ws.Range("N138").Value = -16840.7

$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 254.76471
$ws.Range("I5").Value = 111.454544
$ws.Range("J5").Value = 517.5
$ws.Range("K5").Value = 111.454544
$ws.Range("L5").Value = 517.5
$ws.Range("M5").Value = 0.5454560000000015
$ws.Range("N5").Value = -741.5
# Row 35
$ws.Range("H35").Value = 3874.75
$ws.Range("J35").Value = 1500
$ws.Range("L35").Value = 1500
$ws.Range("N35").Value = -2312
# Row 88
$ws.Range("H88").Value = 2658.516
$ws.Range("I88").Value = 2315
$ws.Range("J88").Value = 2740.96
$ws.Range("K88").Value = 2315
$ws.Range("L88").Value = 2740.96
$ws.Range("M88").Value = -1909
$ws.Range("N88").Value = -3552.96
# Row 91
$ws.Range("H91").Value = 2658.516
$ws.Range("I91").Value = 2315
$ws.Range("J91").Value = 2740.96
$ws.Range("K91").Value = 2315
$ws.Range("L91").Value = 2740.96
$ws.Range("M91").Value = -911
$ws.Range("N91").Value = -5548.96

$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 254.76471
$ws.Range("I4").Value = 111.454544
$ws.Range("J4").Value = 517.5
$ws.Range("K4").Value = 111.454544
$ws.Range("L4").Value = 517.5
$ws.Range("M4").Value = 3.545456000000001
$ws.Range("N4").Value = -747.5
# Row 54
$ws.Range("H54").Value = 22798.4
$ws.Range("I54").Value = 22798.4
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 22798.4
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -22314.4
$ws.Range("N54").ClearContents()
# Row 86
$ws.Range("H86").Value = 1857.2727
$ws.Range("I86").Value = 2061.889
$ws.Range("J86").Value = 936.5
$ws.Range("K86").Value = 2061.889
$ws.Range("L86").Value = 936.5
$ws.Range("M86").Value = -938.8890000000001
$ws.Range("N86").Value = -3182.5
# Row 89
$ws.Range("H89").Value = 1857.2727
$ws.Range("I89").Value = 2061.889
$ws.Range("J89").Value = 936.5
$ws.Range("K89").Value = 10309.445
$ws.Range("L89").Value = 4682.5
$ws.Range("M89").Value = -4693.445
$ws.Range("N89").Value = -15914.5
# Row 107
$ws.Range("H107").Value = 2961.7856
$ws.Range("J107").Value = 1515.6666
$ws.Range("L107").Value = 1515.6666
$ws.Range("N107").Value = -5355.6666

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 3647.45
$ws.Range("I58").Value = 2236.5557
$ws.Range("K58").Value = 2236.5557
$ws.Range("M58").Value = -2033.5557
# Row 62
$ws.Range("H62").Value = 595817.25
$ws.Range("I62").Value = 840967.25
$ws.Range("J62").Value = 7457.2
$ws.Range("K62").Value = 840967.25
$ws.Range("L62").Value = 7457.2
$ws.Range("M62").Value = -840343.25
$ws.Range("N62").Value = -8705.200000000001
# Row 65
$ws.Range("H65").Value = 595817.25
$ws.Range("I65").Value = 840967.25
$ws.Range("J65").Value = 7457.2
$ws.Range("K65").Value = 4204836.25
$ws.Range("L65").Value = 37286
$ws.Range("M65").Value = -4201716.25
$ws.Range("N65").Value = -43526
# Row 136
$ws.Range("H136").Value = 3647.45
$ws.Range("I136").Value = 2236.5557
$ws.Range("K136").Value = 6709.6671
$ws.Range("M136").Value = -4159.6671
# Row 141
$ws.Range("H141").Value = 311077.38
$ws.Range("J141").Value = 311077.38
$ws.Range("L141").Value = 311077.38
$ws.Range("N141").Value = -321437.38

$ws = $wb.Worksheets.Item("CUL")
# Row 9
$ws.Range("H9").Value = 115310.836
$ws.Range("J9").Value = 136373.2
$ws.Range("L9").Value = 409119.6
$ws.Range("N9").Value = -409567.6
# Row 56
$ws.Range("H56").Value = 5925
$ws.Range("I56").Value = 5925
$ws.Range("K56").Value = 5925
$ws.Range("M56").Value = -5395
# Row 97
$ws.Range("H97").Value = 8930247
$ws.Range("J97").Value = 2234.6667
$ws.Range("L97").Value = 6704.000100000001
$ws.Range("N97").Value = -7696.000100000001
# Row 139
$ws.Range("H139").Value = 3605.923
$ws.Range("J139").Value = 5000
$ws.Range("L139").Value = 15000
$ws.Range("N139").Value = -25280

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 160.95
$ws.Range("I2").Value = 59.11111
$ws.Range("J2").Value = 244.27272
$ws.Range("K2").Value = 59.11111
$ws.Range("L2").Value = 244.27272
$ws.Range("M2").Value = 53.88889
$ws.Range("N2").Value = -470.27272
# Row 6
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").ClearContents()
# Row 16
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
# Row 80
$ws.Range("H80").Value = 3430.611
$ws.Range("I80").Value = 3201.75
$ws.Range("K80").Value = 3201.75
$ws.Range("M80").Value = -2203.75
# Row 83
$ws.Range("H83").Value = 3430.611
$ws.Range("I83").Value = 3201.75
$ws.Range("K83").Value = 16008.75
$ws.Range("M83").Value = -11016.75
# Row 113
$ws.Range("H113").Value = 3454.65
$ws.Range("I113").Value = 2259
$ws.Range("J113").Value = 4650.3
$ws.Range("K113").Value = 2259
$ws.Range("L113").Value = 4650.3
$ws.Range("M113").Value = -89
$ws.Range("N113").Value = -8990.299999999999

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 4151.36
$ws.Range("I40").Value = 3567.5789
$ws.Range("K40").Value = 3567.5789
$ws.Range("M40").Value = -3431.5789
# Row 82
$ws.Range("H82").Value = 611.75
$ws.Range("I82").Value = 599
$ws.Range("J82").Value = 624.5
$ws.Range("K82").Value = 599
$ws.Range("L82").Value = 624.5
$ws.Range("M82").Value = -238
$ws.Range("N82").Value = -1346.5
# Row 85
$ws.Range("H85").Value = 611.75
$ws.Range("I85").Value = 599
$ws.Range("J85").Value = 624.5
$ws.Range("K85").Value = 599
$ws.Range("L85").Value = 624.5
$ws.Range("M85").Value = 649
$ws.Range("N85").Value = -3120.5

$ws = $wb.Worksheets.Item("WVR")
# Row 52
$ws.Range("H52").Value = 30245
$ws.Range("I52").Value = 19995
$ws.Range("J52").Value = 40495
$ws.Range("K52").Value = 19995
$ws.Range("L52").Value = 40495
$ws.Range("M52").Value = -19769
$ws.Range("N52").Value = -40947
# Row 100
$ws.Range("H100").Value = 1562.4546
$ws.Range("I100").Value = 2297.4
$ws.Range("J100").Value = 950
$ws.Range("K100").Value = 4594.8
$ws.Range("L100").Value = 1900
$ws.Range("M100").Value = -4053.8
$ws.Range("N100").Value = -2982
